$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the four paragraphs of interest by scanning for unique text.
# (Paragraph indices are stable in the source document, but we scan
#  defensively so the script degrades gracefully if anything shifts.)
# ------------------------------------------------------------------
$iQc        = -1   # "...[QC - this is done]"
$iOliver    = -1   # "...bifurcation time or the bifurcation direction]"
$iThird     = -1   # "3. As the reviewer suggested..."
$iFourth    = -1   # "4. My student and I are investigating..."

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($iQc -eq -1 -and $t.Contains("QC - this is done")) { $iQc = $i }
    if ($iOliver -eq -1 -and $t.Contains("bifurcation time or the bifurcation direction")) { $iOliver = $i }
    if ($iThird -eq -1 -and $t.Contains("3. As the reviewer suggested")) { $iThird = $i }
    if ($iFourth -eq -1 -and $t.Contains("4. My student and I are investigating")) { $iFourth = $i }
}

# ------------------------------------------------------------------
# 1) "[QC - this is done]" -> "[QC - this is done, see updated Table 4]"
# ------------------------------------------------------------------
$d.Paragraphs($iQc).Range.Find.Execute(
    "[QC - this is done]", $false, $false, $false, $false, $false,
    $true, 1, $false,
    "[QC - this is done, see updated Table 4]", 2) | Out-Null

# ------------------------------------------------------------------
# 2) "...bifurcation time or the bifurcation direction]" ->
#    "...bifurcation time nor the bifurcation direction]"
# ------------------------------------------------------------------
$d.Paragraphs($iOliver).Range.Find.Execute(
    "bifurcation time or the bifurcation direction",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "bifurcation time nor the bifurcation direction", 2) | Out-Null

# ------------------------------------------------------------------
# 3) Paragraph 3 text expanded:
#    "...approach, include it in the "Conclusions" section." ->
#    "...approach, and include it in the "Conclusions" section or add
#     a remark after the numerical example section."
# ------------------------------------------------------------------
$p3 = $d.Paragraphs($iThird)
$p3.Range.Find.Execute(
    "single/double stage approach, include it in the",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "single/double stage approach, and include it in the", 2) | Out-Null

$p3 = $d.Paragraphs($iThird)
$p3.Range.Find.Execute(
    [char]8220 + "Conclusions" + [char]8221 + " section.",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    [char]8220 + "Conclusions" + [char]8221 + " section or add a remark after the numerical example section.", 2) | Out-Null

# ------------------------------------------------------------------
# 4) Insert a brand-new paragraph right after paragraph 3, carrying
#    the same paragraph/character formatting (indent, justification,
#    italic, blue font).
# ------------------------------------------------------------------
$p3 = $d.Paragraphs($iThird)
$p3.Range.InsertParagraphAfter() | Out-Null

$iNew = $iThird + 1
$newPara = $d.Paragraphs($iNew)
$newPara.Range.Text = "4. The robustness of the algorithm test (using a single initial point) is essentially eliminating the sampling stage."

# ------------------------------------------------------------------
# 5) Relocate the "_GoBack" bookmark so it again sits at the end of
#    paragraph 3 (its last edit location), rather than where it used
#    to live (inside what is now the Oliver paragraph).
# ------------------------------------------------------------------
try {
    $old = $d.Bookmarks("_GoBack")
    $old.Delete()
} catch {
}

$p3 = $d.Paragraphs($iThird)
$endPoint = $p3.Range.Duplicate
$endPoint.Collapse(0)
try {
    $d.Bookmarks.Add("_GoBack", $endPoint) | Out-Null
} catch {
}

# ------------------------------------------------------------------
# 6) The old "4. My student..." paragraph becomes "5. My student...".
#    It is now two paragraphs further down than $iThird (the expanded
#    3rd paragraph, then the newly-inserted 4th paragraph).
# ------------------------------------------------------------------
$iFifth = $iNew + 1
$p5 = $d.Paragraphs($iFifth)
$p5.Range.Find.Execute(
    "4", $false, $false, $false, $false, $false,
    $true, 1, $false, "5", 2) | Out-Null

Write-Host "Done. iQc=$iQc iOliver=$iOliver iThird=$iThird iNew=$iNew iFifth=$iFifth"
